$wb = $excel.ActiveWorkbook

# --- Sheet1: update price input and its number format, update view ---
$ws1 = $wb.Worksheets.Item("Sheet1")
$ws1.Select()
$ws1.Range("D2").Value = 66
$ws1.Range("D2").NumberFormat = "#,##0.00"

$excel.ActiveWindow.Zoom = 220
$ws1.Range("A2").Select()

# --- Sheet2: discount-rate input, new forecast years, view changes ---
$ws2 = $wb.Worksheets.Item("Sheet2")
$ws2.Select()

# Discount rate input drives the big recalculation ripple (row 12, row 32, etc.)
$ws2.Range("W18").Value = 0.04

# Re-format the per-share NPV cell to two decimals (W22)
$ws2.Range("W22").NumberFormat = "#,##0.00"

# Extend the year header row with five more years (2030-2034)
$ws2.Range("U1").Formula = "=T1+1"
$ws2.Range("U1").NumberFormat = $ws2.Range("T1").NumberFormat
$ws2.Range("V1").Formula = "=U1+1"
$ws2.Range("V1").NumberFormat = $ws2.Range("T1").NumberFormat
$ws2.Range("W1").Formula = "=V1+1"
$ws2.Range("W1").NumberFormat = $ws2.Range("T1").NumberFormat
$ws2.Range("X1").Formula = "=W1+1"
$ws2.Range("X1").NumberFormat = $ws2.Range("T1").NumberFormat
$ws2.Range("Y1").Formula = "=X1+1"
$ws2.Range("Y1").NumberFormat = $ws2.Range("T1").NumberFormat

# Restore the pane/selection state on Sheet2 (frozen header row/col)
$ws2.Range("T1").Select()
